$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.073280345737345
$ws.Range("D2").Value = 1.074935624116435
$ws.Range("E2").Value = 1.073788146056183
$ws.Range("F2").Value = 1.085898064434198
$ws.Range("I2").Value = 1.062838976369312
$ws.Range("J2").Value = 1.078195648215931
$ws.Range("K2").Value = 1.077623521230847
$ws.Range("L2").Value = 1.076479079362058
$ws.Range("M2").Value = 1.088557304461367
$ws.Range("N2").Value = 1.079726809707167
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.074490341832935
$ws.Range("D3").Value = 1.07589904721996
$ws.Range("E3").Value = 1.07487097545955
$ws.Range("F3").Value = 1.086997221257491
$ws.Range("I3").Value = 1.06327727759803
$ws.Range("J3").Value = 1.079063045818735
$ws.Range("K3").Value = 1.078403805076103
$ws.Range("L3").Value = 1.077378256750572
$ws.Range("M3").Value = 1.089475064848211
$ws.Range("N3").Value = 1.080595439114059
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.075272999948233
$ws.Range("D4").Value = 1.076522159017157
$ws.Range("E4").Value = 1.075571876483327
$ws.Range("F4").Value = 1.087708569241472
$ws.Range("I4").Value = 1.063559487409291
$ws.Range("J4").Value = 1.079623447150523
$ws.Range("K4").Value = 1.078907789216614
$ws.Range("L4").Value = 1.077959721671917
$ws.Range("M4").Value = 1.090068437020389
$ws.Range("N4").Value = 1.081156636280017
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.075601961700283
$ws.Range("D5").Value = 1.076784047168979
$ws.Range("E5").Value = 1.075866592942895
$ws.Range("F5").Value = 1.088007649484622
$ws.Range("I5").Value = 1.063677793741126
$ws.Range("J5").Value = 1.079858834325368
$ws.Range("K5").Value = 1.079119446717996
$ws.Range("L5").Value = 1.078204083313955
$ws.Range("M5").Value = 1.090317776529258
$ws.Range("N5").Value = 1.081392357731649
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.075657191906378
$ws.Range("D6").Value = 1.076828015373486
$ws.Range("E6").Value = 1.075916080535221
$ws.Range("F6").Value = 1.08805786813012
$ws.Range("I6").Value = 1.063697638286435
$ws.Range("J6").Value = 1.07989834484475
$ws.Range("K6").Value = 1.079154972217615
$ws.Range("L6").Value = 1.078245107673444
$ws.Range("M6").Value = 1.0903596350593
$ws.Range("N6").Value = 1.081431924360501
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.075277395815132
$ws.Range("D7").Value = 1.076525658645901
$ws.Range("E7").Value = 1.07557581427161
$ws.Range("F7").Value = 1.087712565449178
$ws.Range("I7").Value = 1.063561069538592
$ws.Range("J7").Value = 1.079626593211824
$ws.Range("K7").Value = 1.078910618247827
$ws.Range("L7").Value = 1.077962987182025
$ws.Range("M7").Value = 1.090071769154647
$ws.Range("N7").Value = 1.081159786809086
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.073689330750928
$ws.Range("D8").Value = 1.075261277419398
$ws.Range("E8").Value = 1.074154044461827
$ws.Range("F8").Value = 1.086269505269966
$ws.Range("I8").Value = 1.062987392633449
$ws.Range("J8").Value = 1.07848896859562
$ws.Range("K8").Value = 1.077887411077961
$ws.Range("L8").Value = 1.076783036139567
$ws.Range("M8").Value = 1.088867565784741
$ws.Range("N8").Value = 1.08002054663544
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.070888654597183
$ws.Range("D9").Value = 1.073031041198305
$ws.Range("E9").Value = 1.071650495751559
$ws.Range("F9").Value = 1.083727531065973
$ws.Range("I9").Value = 1.061965748975837
$ws.Range("J9").Value = 1.076477672914076
$ws.Range("K9").Value = 1.076077369545346
$ws.Range("L9").Value = 1.074701000676719
$ws.Range("M9").Value = 1.08674189189661
$ws.Range("N9").Value = 1.078006394683329
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.069019875583096
$ws.Range("D10").Value = 1.0715426642321
$ws.Range("E10").Value = 1.069982621152986
$ws.Range("F10").Value = 1.082033420282699
$ws.Range("I10").Value = 1.061277386202664
$ws.Range("J10").Value = 1.075132266516424
$ws.Range("K10").Value = 1.074865901146401
$ws.Range("L10").Value = 1.07331104146077
$ws.Range("M10").Value = 1.085322226125674
$ws.Range("N10").Value = 1.076659077654271
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.068210253567282
$ws.Range("D11").Value = 1.070897800285436
$ws.Range("E11").Value = 1.069260676486131
$ws.Range("F11").Value = 1.081299966406328
$ws.Range("I11").Value = 1.060977584931758
$ws.Range("J11").Value = 1.074548597497517
$ws.Range("K11").Value = 1.07434017684681
$ws.Range("L11").Value = 1.072708703800478
$ws.Range("M11").Value = 1.084706877844398
$ws.Range("N11").Value = 1.07607457975841
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.067909457208177
$ws.Range("D12").Value = 1.070658209842969
$ws.Range("E12").Value = 1.068992551496848
$ws.Range("F12").Value = 1.081027543875997
$ws.Range("I12").Value = 1.060865963786448
$ws.Range("J12").Value = 1.074331630245602
$ws.Range("K12").Value = 1.074144725393502
$ws.Range("L12").Value = 1.072484896207446
$ws.Range("M12").Value = 1.084478215222693
$ws.Range("N12").Value = 1.075857304388111
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.06797398210387
$ws.Range("D13").Value = 1.07070960550606
$ws.Range("E13").Value = 1.06905006354166
$ws.Range("F13").Value = 1.081085978779037
$ws.Range("I13").Value = 1.060889918752727
$ws.Range("J13").Value = 1.074378178019565
$ws.Range("K13").Value = 1.074186658288585
$ws.Range("L13").Value = 1.072532906996252
$ws.Range("M13").Value = 1.084527268443333
$ws.Range("N13").Value = 1.075903918265252
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.068185391020462
$ws.Range("D14").Value = 1.070877996880172
$ws.Range("E14").Value = 1.069238512433981
$ws.Range("F14").Value = 1.08127744757817
$ws.Range("I14").Value = 1.060968363634598
$ws.Range("J14").Value = 1.074530666323592
$ws.Range("K14").Value = 1.074324024325986
$ws.Range("L14").Value = 1.072690205281327
$ws.Range("M14").Value = 1.084687978462926
$ws.Range("N14").Value = 1.076056623120162
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.068315638123826
$ws.Range("D15").Value = 1.070981740476901
$ws.Range("E15").Value = 1.069354626940812
$ws.Range("F15").Value = 1.08139541976197
$ws.Range("I15").Value = 1.061016661424951
$ws.Range("J15").Value = 1.074624597282125
$ws.Range("K15").Value = 1.074408636977516
$ws.Range("L15").Value = 1.072787112299598
$ws.Range("M15").Value = 1.084786984603432
$ws.Range("N15").Value = 1.07615068747143
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.069073598305234
$ws.Range("D16").Value = 1.071585453480288
$ws.Range("E16").Value = 1.070030539546656
$ws.Range("F16").Value = 1.082082099356916
$ws.Range("I16").Value = 1.061297246385571
$ws.Range("J16").Value = 1.075170979424749
$ws.Range("K16").Value = 1.074900767395499
$ws.Range("L16").Value = 1.073351006507641
$ws.Range("M16").Value = 1.085363051524759
$ws.Range("N16").Value = 1.076697845539368
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.069548930365839
$ws.Range("D17").Value = 1.071964042314244
$ws.Range("E17").Value = 1.070454589465691
$ws.Range("F17").Value = 1.082512863183551
$ws.Range("I17").Value = 1.061472784640686
$ws.Range("J17").Value = 1.075513415203127
$ws.Range("K17").Value = 1.07520915894438
$ws.Range("L17").Value = 1.073704594239788
$ws.Range("M17").Value = 1.085724235288954
$ws.Range("N17").Value = 1.077040767615828
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.069826142342382
$ws.Range("D18").Value = 1.072184829523635
$ws.Range("E18").Value = 1.070701955386387
$ws.Range("F18").Value = 1.082764130869688
$ws.Range("I18").Value = 1.061575005727792
$ws.Range("J18").Value = 1.075713046244337
$ws.Range("K18").Value = 1.075388927593444
$ws.Range("L18").Value = 1.073910790111137
$ws.Range("M18").Value = 1.085934847576948
$ws.Range("N18").Value = 1.077240682156015
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.069920657549173
$ws.Range("D19").Value = 1.072260105971529
$ws.Range("E19").Value = 1.070786305007134
$ws.Range("F19").Value = 1.082849808463487
$ws.Range("I19").Value = 1.061609832106947
$ws.Range("J19").Value = 1.075781097342196
$ws.Range("K19").Value = 1.075450205278201
$ws.Range("L19").Value = 1.073981089825003
$ws.Range("M19").Value = 1.08600665076715
$ws.Range("N19").Value = 1.077308829894238
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.069497936010976
$ws.Range("D20").Value = 1.071923427177554
$ws.Range("E20").Value = 1.070409090368348
$ws.Range("F20").Value = 1.082466645232226
$ws.Range("I20").Value = 1.061453968369553
$ws.Range("J20").Value = 1.075476686055687
$ws.Range("K20").Value = 1.075176082932582
$ws.Range("L20").Value = 1.073666662375366
$ws.Range("M20").Value = 1.085685489914151
$ws.Range("N20").Value = 1.077003986308786
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.068123138223611
$ws.Range("D21").Value = 1.070828411462526
$ws.Range("E21").Value = 1.069183017908147
$ws.Range("F21").Value = 1.081221064387657
$ws.Range("I21").Value = 1.060945270800417
$ws.Range("J21").Value = 1.074485766929115
$ws.Range("K21").Value = 1.074283578265689
$ws.Range("L21").Value = 1.072643886868414
$ws.Range("M21").Value = 1.084640655996136
$ws.Range("N21").Value = 1.076011659963395
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.067258361259697
$ws.Range("D22").Value = 1.070139588197562
$ws.Range("E22").Value = 1.068412353374405
$ws.Range("F22").Value = 1.08043800385395
$ws.Range("I22").Value = 1.060623918467214
$ws.Range("J22").Value = 1.073861772183345
$ws.Range("K22").Value = 1.073721417621505
$ws.Range("L22").Value = 1.072000406952975
$ws.Range("M22").Value = 1.083983177960259
$ws.Range("N22").Value = 1.075386779073512
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.067716833357307
$ws.Range("D23").Value = 1.070504779375188
$ws.Range("E23").Value = 1.068820877117011
$ws.Range("F23").Value = 1.080853111354995
$ws.Range("I23").Value = 1.060794417208873
$ws.Range("J23").Value = 1.074192655519447
$ws.Range("K23").Value = 1.074019525517096
$ws.Range("L23").Value = 1.072341568020987
$ws.Range("M23").Value = 1.08433177198534
$ws.Range("N23").Value = 1.075718132301905
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.069520978292172
$ws.Range("D24").Value = 1.071941779524793
$ws.Range("E24").Value = 1.070429649375549
$ws.Range("F24").Value = 1.082487529102903
$ws.Range("I24").Value = 1.061462471149845
$ws.Range("J24").Value = 1.075493282703881
$ws.Range("K24").Value = 1.075191028901389
$ws.Range("L24").Value = 1.073683802293995
$ws.Range("M24").Value = 1.085702997467031
$ws.Range("N24").Value = 1.077020606526124
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.071612981060148
$ws.Range("D25").Value = 1.073607880727753
$ws.Range("E25").Value = 1.072297515607506
$ws.Range("F25").Value = 1.084384593382809
$ws.Range("I25").Value = 1.062231146408685
$ws.Range("J25").Value = 1.076998436759596
$ws.Range("K25").Value = 1.076546146137173
$ws.Range("L25").Value = 1.075239594563381
$ws.Range("M25").Value = 1.087291875089696
$ws.Range("N25").Value = 1.078527898073242
